$wb = $excel.ActiveWorkbook

# Rename first sheet from "WWR per cluster" to "Baseline"
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Baseline"

# Update the saved selection on the (now) "Baseline" sheet from J19 to B21
$ws.Range("B21").Select()
